$d = $word.ActiveDocument

# Replace the signatory name "LEONARDO SILVERIO FERREIRA" with
# "MANOEL JEFETE DA SILVA TENONIO" wherever it appears in the document body.
$d.Content.Find.Execute(
    "LEONARDO SILVERIO FERREIRA",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "MANOEL JEFETE DA SILVA TENONIO",
    2
)
